$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target row order per year-block: Oct, Nov, Dec, Jan, Feb, ..., Sep
# (A: period label, B/C/D: the three series values)
$data = @(
    @("2014-10", 101.12, 100.4743, 98.9387),
    @("2014-11", 101.1257, 100.649, 99.5475),
    @("2014-12", 100.6284, 100.6566, 99.85290000000001),
    @("2014-01", 100.7959, 99.79649999999999, 99.9545),
    @("2014-02", 101.1775, 100.5356, 99.962),
    @("2014-03", 100.6482, 100.6973, 100.0345),
    @("2014-04", 101.0813, 100.6977, 99.9541),
    @("2014-05", 101.0066, 100.4889, 99.6015),
    @("2014-06", 101.1092, 100.4829, 99.53060000000001),
    @("2014-07", 101.2366, 100.5859, 99.53230000000001),
    @("2014-08", 101.2756, 100.5901, 99.229),
    @("2014-09", 101.4698, 100.644, 98.8681),
    @("2015-10", 99.7, 98.09999999999999, 98.40000000000001),
    @("2015-11", 98.786, 98.05800000000001, 98.40389999999999),
    @("2015-12", 97.2306, 98.1395, 97.9825),
    @("2015-01", 100.6231, 100.6015, 99.6966),
    @("2015-02", 100.4987, 99.8319, 99.6354),
    @("2015-03", 100.4911, 99.56610000000001, 99.54940000000001),
    @("2015-04", 100.6833, 99.3754, 99.52889999999999),
    @("2015-05", 100.5111, 98.5341, 98.6245),
    @("2015-06", 100.4419, 98.2679, 98.09780000000001),
    @("2015-07", 100.2568, 97.9821, 97.8622),
    @("2015-08", 100.1791, 97.9888, 97.7747),
    @("2015-09", 99.8785, 97.9532, 97.9958),
    @("2016-10", 97.5, 98.40000000000001, 94.3),
    @("2016-11", 98.09999999999999, 98.3, 93.90000000000001),
    @("2016-12", 99.5, 98.09999999999999, 94.09999999999999),
    @("2016-01", 97.5397, 97.3724, 95.8382),
    @("2016-02", 97.3211, 97.30500000000001, 94.2354),
    @("2016-03", 96.8725, 97.4862, 93.72150000000001),
    @("2016-04", 96.67910000000001, 97.4481, 93.30840000000001),
    @("2016-05", 96.8, 98.3, 94.40000000000001),
    @("2016-06", 97, 98.7, 95.2),
    @("2016-07", 97, 98.8, 95),
    @("2016-08", 97, 98.7, 94.90000000000001),
    @("2016-09", 97.2, 98.59999999999999, 94.7),
    @("2017-10", 102.5, 98.8, 101.1),
    @("2017-11", 102.3, 98.90000000000001, 101.2),
    @("2017-12", 102.6, 98.90000000000001, 101),
    @("2017-01", 100.1, 98.7, 96.3),
    @("2017-02", 101, 99, 97.90000000000001),
    @("2017-03", 101.4, 99, 98.2),
    @("2017-04", 102, 99.2, 98.7),
    @("2017-05", 102.1, 99.59999999999999, 98.8),
    @("2017-06", 102.1, 99.7, 99),
    @("2017-07", 102.5, 99.90000000000001, 99.59999999999999),
    @("2017-08", 102.4, 99.59999999999999, 100.3),
    @("2017-09", 102.4, 98.8, 100.8)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $data[$i]
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
}
